$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, E are plain text/URLs/percentages that Excel will not
# reinterpret as numbers, so they can be assigned directly.
# Column D contains numeric-looking strings (prices) that must remain
# text, matching the original inline-string cells. We force text entry
# by switching the cell to the "@" (Text) number format before the
# assignment, then clear the format again so the cell keeps its
# original (default) style once the text value is safely stored.

$ws.Range("E2").Value = "  -6.78%  "
$ws.Range("E3").Value = "  -8.75%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("E5").Value = "  -10.88%  "
$ws.Range("E6").Value = "  -15.83%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -9.34%  "
$ws.Range("E9").Value = "  -17.10%  "
$ws.Range("E10").Value = "  -16.52%  "
$ws.Range("E11").Value = "  -8.74%  "
$ws.Range("E12").Value = "  -13.26%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("E13").Value = "  -16.95%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("E14").Value = "  -20.43%  "
$ws.Range("E15").Value = "  -9.61%  "
$ws.Range("E16").Value = "  -7.53%  "
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("E18").Value = "  -9.64%  "
$ws.Range("E19").Value = "  -10.55%  "
$ws.Range("E20").Value = "  -14.38%  "
$ws.Range("E21").Value = "  -16.85%  "
$ws.Range("E22").Value = "  -16.67%  "
$ws.Range("E23").Value = "  -19.42%  "
$ws.Range("E24").Value = "  -12.86%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E25").Value = "  -13.76%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  -18.61%  "
$ws.Range("E28").Value = "  -14.15%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E29").Value = "  -17.42%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E30").Value = "  -15.86%  "
$ws.Range("B31").Value = "Mantle"
$ws.Range("C31").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E31").Value = "  -8.18%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("E33").Value = "  -14.53%  "
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("E35").Value = "  -12.23%  "
$ws.Range("E36").Value = "  -15.96%  "
$ws.Range("E37").Value = "  -19.10%  "
$ws.Range("E38").Value = "  -11.74%  "
$ws.Range("E39").Value = "  -7.73%  "
$ws.Range("E40").Value = "  -15.19%  "
$ws.Range("E41").Value = "  -17.81%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("E42").Value = "  -11.33%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  -17.15%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("E45").Value = "  -18.26%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E46").Value = "  -7.16%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E47").Value = "  -14.37%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("E49").Value = "  -20.98%  "
$ws.Range("B50").Value = "BitgetToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("E50").Value = "  -4.45%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E51").Value = "  -19.14%  "

$dCells = @("D2", "D3", "D4", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.663.86"
$ws.Range("D3").Value = "2.955.63"
$ws.Range("D4").Value = "0.995"
$ws.Range("D5").Value = "528.81"
$ws.Range("D6").Value = "128.90"
$ws.Range("D8").Value = "2.929.41"
$ws.Range("D9").Value = "0.451"
$ws.Range("D10").Value = "0.144"
$ws.Range("D11").Value = "5.93"
$ws.Range("D12").Value = "0.427"
$ws.Range("D13").Value = "0.0000201"
$ws.Range("D14").Value = "31.02"
$ws.Range("D15").Value = "3.411.21"
$ws.Range("D16").Value = "62.176.62"
$ws.Range("D17").Value = "0.110"
$ws.Range("D18").Value = "2.927.43"
$ws.Range("D19").Value = "476.48"
$ws.Range("D20").Value = "6.06"
$ws.Range("D21").Value = "12.36"
$ws.Range("D22").Value = "0.631"
$ws.Range("D23").Value = "6.37"
$ws.Range("D24").Value = "74.50"
$ws.Range("D25").Value = "11.82"
$ws.Range("D26").Value = "0.998"
$ws.Range("D27").Value = "2.60"
$ws.Range("D28").Value = "6.92"
$ws.Range("D29").Value = "1.80"
$ws.Range("D30").Value = "24.45"
$ws.Range("D31").Value = "1.06"
$ws.Range("D32").Value = "0.991"
$ws.Range("D33").Value = "2.28"
$ws.Range("D34").Value = "51.43"
$ws.Range("D35").Value = "467.95"
$ws.Range("D36").Value = "5.40"
$ws.Range("D37").Value = "4.65"
$ws.Range("D38").Value = "0.0376"
$ws.Range("D39").Value = "0.113"
$ws.Range("D40").Value = "0.0727"
$ws.Range("D41").Value = "7.63"
$ws.Range("D42").Value = "2.601.04"
$ws.Range("D43").Value = "0.999"
$ws.Range("D44").Value = "2.27"
$ws.Range("D45").Value = "0.216"
$ws.Range("D46").Value = "109.97"
$ws.Range("D47").Value = "0.0981"
$ws.Range("D48").Value = "1.80"
$ws.Range("D49").Value = "0.0₃0465"
$ws.Range("D50").Value = "1.20"
$ws.Range("D51").Value = "21.32"

foreach ($addr in $dCells) {
    $ws.Range($addr).ClearFormats()
}

